$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.196.21'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.37%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.328.61'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.61%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '530.63'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.23'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.58%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.993'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.45%  '
$ws.Range('E8').Value = '  -0.69%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.353.38'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('E10').Value = '  -1.40%  '
$ws.Range('E11').Value = '  +0.49%  '
$ws.Range('E12').Value = '  -2.25%  '
$ws.Range('E13').Value = '  +0.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.746.88'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.55'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '57.211.67'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.43%  '
$ws.Range('E17').Value = '  -1.30%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.348.48'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '337.44'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.44'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.92'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.58%  '
$ws.Range('E22').Value = '  -1.65%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.77'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.87'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +10.93%  '
$ws.Range('E26').Value = '  -0.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.993'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.34'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '169.79'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('E30').Value = '  +1.23%  '
$ws.Range('E31').Value = '  -2.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.12'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.52'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.991'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.26'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.88%  '
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.911'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('E39').Value = '  +1.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '38.97'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '148.11'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.77%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.378'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.59'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.72%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '284.58'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.44%  '
$ws.Range('E45').Value = '  -1.43%  '
$ws.Range('E46').Value = '  -0.37%  '
$ws.Range('E47').Value = '  -0.62%  '
$ws.Range('E48').Value = '  -0.29%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.74'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.39%  '
$ws.Range('E50').Value = '  -0.98%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.34'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.68%  '
